# PRLB Quarterly Financials update
# Inserts two new quarterly columns (FY2018-Q4 ending 2018-12-31 and
# FY2018-Q3 ending 2018-09-30) ahead of the existing data, shifting the
# previously-existing quarters two columns to the right, then fills in
# the new quarter figures reported for the three statements
# (Income Statement, Balance Sheet, Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank columns at D:E - this shifts the old D:K data (and
#    its number formats / cell styles) to F:M automatically.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2. Carry over the (now shifted) column formatting from F:G into the
#    new D:E columns so the new cells use the same styles (date format
#    for the header rows, number format for the data rows) without
#    introducing any new style entries. This is only done for the rows
#    that belong to the three statement tables (the blank spacer/title
#    rows in between never had data in D:K and must stay untouched).
$dataRowRanges = @(
    @(7, 35),
    @(38, 77),
    @(80, 102)
)
foreach ($rng in $dataRowRanges) {
    $r1 = $rng[0]
    $r2 = $rng[1]
    $ws.Range("F$r1`:G$r2").Copy()
    $ws.Range("D$r1`:E$r2").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# 3. Populate the new D (period ending 2018-12-31) and E (period ending
#    2018-09-30) columns with the reported values for each row.
$rows = @(
    @(7, 43465, 43373),
    @(8, 112800, 115400),
    @(9, 53600, 53000),
    @(10, 59200, 62400),
    @(12, 7600, 7500),
    @(13, 0, 0),
    @(14, 0, 0),
    @(15, 0, 0),
    @(17, 92600, 90400),
    @(18, 20200, 25000),
    @(20, 1400, 400),
    @(21, 28800, 32200),
    @(22, 0, 0),
    @(23, 21500, 25400),
    @(24, -1900, 4000),
    @(25, 0, 0),
    @(26, 23400, 21400),
    @(27, 23400, 21400),
    @(28, 0, 0),
    @(29, -4100, -500),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, -1400, -400),
    @(33, 19300, 20900),
    @(34, 0, 0),
    @(35, 19300, 20900),
    @(38, 43465, 43373),
    @(41, 85000, 62600),
    @(42, 46800, 57600),
    @(43, 64900, 62100),
    @(44, 10100, 9600),
    @(45, 8600, 8600),
    @(46, 215400, 200500),
    @(47, "NA", "NA"),
    @(48, 228000, 210100),
    @(49, 148600, 148500),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 27000, 36800),
    @(53, 0, 0),
    @(54, 619000, 595900),
    @(57, 17400, 17300),
    @(58, 0, 0),
    @(59, 35300, 32100),
    @(60, 52700, 49400),
    @(61, 0, 0),
    @(62, 24800, 15300),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 77500, 64700),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, 291500, 283500),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 541500, 531200),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, 19300, 20900),
    @(83, 7200, 6800),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 38000, 25600),
    @(91, -25200, -18700),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -4800, -27400),
    @(96, 0, 0),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, -10700, 1900),
    @(101, 0, -200),
    @(102, 22500, 0)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    $ws.Cells.Item($r, 4).Value2 = $dVal
    $ws.Cells.Item($r, 5).Value2 = $eVal
}
